$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be written as text (matching the workbook's
    # original inline-string cell type) instead of letting Excel
    # auto-convert numeric-looking strings into numbers, then restore
    # the default "Normal" style so no stray formatting is introduced.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '51.073.09'
$ws.Range('E2').Value = '  -1.27%  '
Set-TextValue 'D3' '2.908.72'
$ws.Range('E3').Value = '  -0.68%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.21%  '
Set-TextValue 'D5' '371.50'
$ws.Range('E5').Value = '  +5.60%  '
Set-TextValue 'D6' '103.86'
$ws.Range('E6').Value = '  -2.91%  '
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -4.01%  '
Set-TextValue 'D10' '36.51'
$ws.Range('E10').Value = '  -3.30%  '
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('E12').Value = '  -1.56%  '
Set-TextValue 'D13' '18.37'
$ws.Range('E13').Value = '  -2.91%  '
Set-TextValue 'D14' '3.367.34'
$ws.Range('E14').Value = '  -0.96%  '
Set-TextValue 'D15' '7.38'
$ws.Range('E15').Value = '  -2.38%  '
Set-TextValue 'D16' '2.909.16'
$ws.Range('E16').Value = '  -1.11%  '
Set-TextValue 'D17' '0.931'
$ws.Range('E17').Value = '  -4.02%  '
Set-TextValue 'D18' '50.972.57'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('E19').Value = '  -7.56%  '
Set-TextValue 'D20' '7.19'
$ws.Range('E20').Value = '  -2.15%  '
Set-TextValue 'D21' '12.84'
$ws.Range('E21').Value = '  -4.25%  '
Set-TextValue 'D22' '0.0₃0941'
$ws.Range('E22').Value = '  -2.23%  '
$ws.Range('E23').Value = '  -1.21%  '
Set-TextValue 'D24' '259.31'
$ws.Range('E24').Value = '  -1.08%  '
Set-TextValue 'D25' '2.70'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('E26').Value = '  +3.59%  '
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('E28').Value = '  +0.01%  '
Set-TextValue 'D29' '25.66'
$ws.Range('E29').Value = '  -3.29%  '
Set-TextValue 'D30' '6.93'
$ws.Range('E30').Value = '  -8.29%  '
$ws.Range('E31').Value = '  -1.09%  '
Set-TextValue 'D32' '6.15'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('E33').Value = '  -3.37%  '
Set-TextValue 'D34' '2.12'
$ws.Range('E34').Value = '  -1.51%  '
Set-TextValue 'D35' '34.57'
$ws.Range('E35').Value = '  -2.99%  '
Set-TextValue 'D36' '50.81'
$ws.Range('E36').Value = '  -0.04%  '
Set-TextValue 'D38' '0.0420'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('E39').Value = '  -3.06%  '
Set-TextValue 'D40' '2.64'
$ws.Range('E40').Value = '  -0.85%  '
Set-TextValue 'D41' '17.03'
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('E42').Value = '  -6.08%  '
$ws.Range('E43').Value = '  -2.28%  '
Set-TextValue 'D44' '22.16'
$ws.Range('E44').Value = '  -1.05%  '
Set-TextValue 'D45' '119.36'
$ws.Range('E45').Value = '  -0.26%  '
Set-TextValue 'D46' '2.09'
$ws.Range('E46').Value = '  -2.45%  '
Set-TextValue 'D47' '2.016.92'
$ws.Range('E47').Value = '  -3.85%  '
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('E49').Value = '  -4.46%  '
$ws.Range('E50').Value = '  +1.68%  '
Set-TextValue 'D51' '0.0309'
$ws.Range('E51').Value = '  -10.12%  '
